$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (must remain plain text, not numeric)
$updates = [ordered]@{
    'D2' = '68.968.06'
    'E2' = '  +1.21%  '
    'D3' = '3.931.27'
    'E3' = '  +0.44%  '
    'D4' = '1.00'
    'E4' = '  +0.17%  '
    'D5' = '487.89'
    'E5' = '  +0.57%  '
    'D6' = '147.67'
    'E6' = '  +1.05%  '
    'E7' = '  -0.37%  '
    'D8' = '0.998'
    'E8' = '  +0.09%  '
    'E9' = '  +1.61%  '
    'D10' = '0.176'
    'E10' = '  +4.64%  '
    'D11' = '0.0000345'
    'E11' = '  -2.78%  '
    'D12' = '43.15'
    'E12' = '  +1.49%  '
    'D13' = '10.50'
    'E13' = '  -0.91%  '
    'D14' = '4.559.15'
    'E14' = '  +0.71%  '
    'D15' = '3.919.50'
    'E15' = '  +0.09%  '
    'D16' = '14.31'
    'E16' = '  -2.81%  '
    'E17' = '  -0.70%  '
    'D18' = '19.98'
    'E18' = '  +0.81%  '
    'E19' = '  +1.83%  '
    'D20' = '69.026.27'
    'E20' = '  +1.20%  '
    'D21' = '436.54'
    'E21' = '  -2.63%  '
    'E22' = '  +4.60%  '
    'D23' = '14.65'
    'E23' = '  -0.20%  '
    'D24' = '89.45'
    'E24' = '  +0.57%  '
    'D25' = '12.18'
    'E25' = '  +10.80%  '
    'E26' = '  +3.44%  '
    'D27' = '11.10'
    'E27' = '  -3.95%  '
    'D28' = '37.17'
    'E28' = '  -4.35%  '
    'E29' = '  -3.24%  '
    'D30' = '710.54'
    'E30' = '  +2.99%  '
    'E31' = '  +1.10%  '
    'D32' = '13.41'
    'E32' = '  +0.18%  '
    'E33' = '  +1.10%  '
    'D34' = '0.478'
    'E34' = '  +30.76%  '
    'D35' = '0.0₃0902'
    'E35' = '  -2.46%  '
    'D36' = '6.09'
    'E36' = '  +6.51%  '
    'D37' = '61.13'
    'E37' = '  +3.46%  '
    'D38' = '40.97'
    'E38' = '  -1.95%  '
    'E39' = '  -0.37%  '
    'D40' = '0.999'
    'E40' = '  +0.07%  '
    'E41' = '  +0.24%  '
    'D42' = '2.96'
    'E42' = '  +1.47%  '
    'E43' = '  +2.28%  '
    'D44' = '3.08'
    'E44' = '  +0.29%  '
    'D45' = '3.02'
    'E45' = '  +0.51%  '
    'D46' = '0.0₆0381'
    'E46' = '  +16.80%  '
    'E47' = '  +9.16%  '
    'E48' = '  +1.11%  '
    'E49' = '  +5.45%  '
    'E50' = '  -1.30%  '
    'D51' = '2.08'
    'E51' = '  -2.06%  '
}

foreach ($addr in $updates.Keys) {
    $text = $updates[$addr]
    $cell = $ws.Range($addr)
    # Build a formula that evaluates to the literal text, so Excel
    # treats it as a string result instead of auto-converting
    # numeric-looking text (e.g. "1.00") into a Double.
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    # Flatten the formula down to a plain stored value (keeps the
    # string type) without leaving a formula or touching styles.
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0
